# Update the Fat4-Dchs1 NATMI sheet with the new TPM-derived values.
# Ligand-side stats are keyed by "Sending cluster" (column A) and drive
# columns E-J; receptor-side stats are keyed by "Target cluster" (column D)
# and drive columns K-P. Columns Q-T are simple products of the ligand and
# receptor columns (avg*avg, total*total, specI*specO, specJ*specP).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New per-sending-cluster ligand stats (E,F,G,H)
$ligand = @{
    "ECs"           = @{ E = 3; F = 1;                  G = 24.00186933333333;  H = 72.005608 }
    "FAPs"          = @{ E = 3; F = 1;                  G = 21.652497;          H = 64.957491 }
    "MuSCs"         = @{ E = 3; F = 1;                  G = 5.091779333333333;  H = 15.275338 }
    "Resolving-Mac" = @{ E = 2; F = 0.6666666666666666; G = 0.04006366666666666;H = 0.120191 }
}

# New per-target-cluster receptor stats (K,L,M,N) - K,L are unchanged
$receptor = @{
    "ECs"           = @{ K = 3; L = 1; M = 13.76681233333333;  N = 41.300437 }
    "FAPs"          = @{ K = 3; L = 1; M = 15.18384366666667;  N = 45.551531 }
    "MuSCs"         = @{ K = 3; L = 1; M = 0.8411156666666667; N = 2.523347 }
    "Resolving-Mac" = @{ K = 3; L = 1; M = 0.2787996666666667; N = 0.8363990000000001 }
}

# Ligand/receptor derived specificity denominators = sum across all clusters
$sumG = 0.0
$sumH = 0.0
foreach ($k in $ligand.Keys) {
    $sumG += $ligand[$k].G
    $sumH += $ligand[$k].H
}

$sumM = 0.0
$sumN = 0.0
foreach ($k in $receptor.Keys) {
    $sumM += $receptor[$k].M
    $sumN += $receptor[$k].N
}

# Data rows are 2..17; column A = sending cluster, column D = target cluster
for ($row = 2; $row -le 17; $row++) {
    $sending = $ws.Cells.Item($row, 1).Value2
    $target  = $ws.Cells.Item($row, 4).Value2

    $lig = $ligand[$sending]
    $rec = $receptor[$target]

    $G = $lig.G
    $H = $lig.H
    $I = $G / $sumG
    $J = $H / $sumH

    $M = $rec.M
    $N = $rec.N
    $O = $M / $sumM
    $P = $N / $sumN

    $ws.Cells.Item($row, 5).Value2  = $lig.E   # E
    $ws.Cells.Item($row, 6).Value2  = $lig.F   # F
    $ws.Cells.Item($row, 7).Value2  = $G       # G
    $ws.Cells.Item($row, 8).Value2  = $H       # H
    $ws.Cells.Item($row, 9).Value2  = $I       # I
    $ws.Cells.Item($row, 10).Value2 = $J       # J

    $ws.Cells.Item($row, 13).Value2 = $M       # M
    $ws.Cells.Item($row, 14).Value2 = $N       # N
    $ws.Cells.Item($row, 15).Value2 = $O       # O
    $ws.Cells.Item($row, 16).Value2 = $P       # P

    $ws.Cells.Item($row, 17).Value2 = $G * $M  # Q
    $ws.Cells.Item($row, 18).Value2 = $H * $N  # R
    $ws.Cells.Item($row, 19).Value2 = $I * $O  # S
    $ws.Cells.Item($row, 20).Value2 = $J * $P  # T
}

Write-Output "updated rows 2-17"
